$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "SOFTWARE"
$ws.Range("C2").Value = "NOV"
$ws.Range("E2").Value = "00005 Shinde Shraddha"
$ws.Range("G2").Value = "01/10/2017"
$ws.Range("H2").Value = "02/10/2017"

# Row 3 updates
$ws.Range("A3").Value = "SOFTWARE"
$ws.Range("C3").Value = "NOV"
$ws.Range("E3").Value = "00005 Shinde Shraddha"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "09"
$ws.Range("G3").Value = "03/10/2017"
$ws.Range("H3").Value = "03/10/2017"

$ws.Range("H3").Select() | Out-Null
